$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J1").EntireColumn.Delete()
$ws.Columns("J:J").Select()

